# Update countries & provincias Spain
# Refresh COVID-19 country data snapshot + update "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: last-updated timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 21:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 551081
$ws.Range("C4").Value = 18202
$ws.Range("D4").Value = 31369
$ws.Range("E4").Value = 498044
$ws.Range("G4").Value = 1091
$ws.Range("H4").Value = 21668

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 127007
$ws.Range("C8").Value = 1555
$ws.Range("E8").Value = 63746
$ws.Range("G8").Value = 90
$ws.Range("H8").Value = 2961

# --- Row 76: Uzbekistan ---
$ws.Range("B76").Value = 865
$ws.Range("C76").Value = 98
$ws.Range("E76").Value = 796

# --- Row 91: Costa Rica ---
$ws.Range("B91").Value = 595
$ws.Range("C91").Value = 18
$ws.Range("D91").Value = 56
$ws.Range("E91").Value = 536
$ws.Range("F91").Value = 14

# --- Row 92: Costa de Marfil ---
$ws.Range("B92").Value = 574
$ws.Range("C92").Value = 41
$ws.Range("D92").Value = 85
$ws.Range("E92").Value = 484
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 5

# --- Rows 143/144: Uganda overtakes Polinesia Francesa in the ranking ---
$ws.Range("A143").Value = "Uganda"
$ws.Range("B143").Value = 54
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 4
$ws.Range("E143").Value = 50
$ws.Range("F143").Value = 0

$ws.Range("A144").Value = "Polinesia Francesa"
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 53
$ws.Range("F144").Value = 1
